# Refresh of the cryptos price/volume snapshot (GitHub Actions scheduled update).
# Price (column D) and Volume(1h) (column E) text values are refreshed in place.
# A few rows also swap rank position with their neighbour (26<->27, 49/50/51 shift
# down one), which shows up as the Coin name/Link/Price/Volume all changing together.
#
# All D/E/B/C cells on this sheet are stored as plain TEXT (t="inlineStr" in the
# source file), not numbers - values like "1.00", "255.83", "98.789.67" (note the
# double '.') or "  +0.32%  " (padded with spaces) must stay text. Excel's COM
# Range.Value setter auto-coerces plain decimal-looking strings to numbers, which
# would corrupt values such as "1.00" -> 1 or "255.83" -> 255.83 (losing its text
# type / exact formatting). To avoid that, column D writes are wrapped with a
# temporary Text number format, then the style is reset back to "Normal" so the
# on-disk cell style is unchanged. Column E's percentage strings keep their
# leading/trailing double-space padding, which already prevents Excel from
# reinterpreting them as numbers, and columns B/C are plain names/URLs, so those
# are set directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '98.803.51'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.34%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.318.09'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.27%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '255.83'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.49%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '632.73'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.52%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.46'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +18.99%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.410'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +6.08%  '
$ws.Range("E9").Value = '  -0.03%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.00'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +22.48%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '3.314.72'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.26%  '
$ws.Range("E12").Value = '  +2.99%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '42.77'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +19.59%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '98.748.14'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.67%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000252'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.00%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.946.65'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.26%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.43'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.36%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.320.08'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.35%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '16.36'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +8.99%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.51'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.89%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.63'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +11.83%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '487.85'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.33%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.61'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.70%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.0000204'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.65%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '6.05'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +5.65%  '
$ws.Range("B26").Value = 'Litecoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '91.20'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.46%  '
$ws.Range("B27").Value = 'Stellar'
$ws.Range("C27").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.339'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +31.61%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '12.35'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.25%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.494.15'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.40%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.147'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +16.94%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '11.22'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +21.44%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.192'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.43%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.998'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.17%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '28.28'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.00%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.488'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +8.51%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '7.49'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.63%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.152'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.45%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.99'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.56%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '504.06'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.46%  '
$ws.Range("E41").Value = '  -0.28%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.87'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.80%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.27'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.14%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.806'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.20%  '
$ws.Range("E45").Value = '  +0.06%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.18'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.31%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.01'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.84%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '160.56'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.51%  '
$ws.Range("B49").Value = 'ImmutableX'
$ws.Range("C49").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.47'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +7.90%  '
$ws.Range("B50").Value = 'Cosmos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.47'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +16.15%  '
$ws.Range("B51").Value = 'Filecoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.84'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +6.56%  '
